{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Reproduces the commit \"ajout des defis rencontrer sur le document word\":\n//   1. The empty paragraph right after the \"CREATION DE LA PAGE \u00ab SITE DORDI \u00bb \u00b8\"\n//      Heading 1 gets a single space \" \" typed into it.\n//   2. The (until now empty) last \"ListParagraph\" bullet (numId 6) at the end of\n//      the document is filled in with the first \"d\u00e9fis rencontr\u00e9s\" bullet, and two\n//      more bullets (same list/style) are appended after it.\n\nconst body = context.document.body;\n\n// --- Change 1: the blank paragraph right after the \"CREATION DE LA PAGE...\" title ---\nconst titleSearch = body.search(\"CREATION DE LA PAGE \u00ab SITE DORDI \u00bb \u00b8\", { matchCase: false });\ntitleSearch.load(\"items\");\nawait context.sync();\n\nif (titleSearch.items.length > 0) {\n  const titlePara = titleSearch.items[0].paragraphs.getFirst();\n  const blankPara = titlePara.getNext();\n  blankPara.load(\"text\");\n  await context.sync();\n  if (!blankPara.text) {\n    blankPara.insertText(\" \", Word.InsertLocation.end);\n  }\n} else {\n  // Fallback: locate the Heading 1 paragraph directly.\n  body.paragraphs.load(\"items/style,items/text\");\n  await context.sync();\n  const titlePara = body.paragraphs.items.find(\n    (p) => p.style === \"Heading 1\" && p.text.indexOf(\"CREATION DE LA PAGE\") !== -1\n  );\n  if (titlePara) {\n    const blankPara = titlePara.getNext();\n    blankPara.insertText(\" \", Word.InsertLocation.end);\n  }\n}\n\n// --- Change 2: fill in / append the \"d\u00e9fis rencontr\u00e9s\" bullet list ---\nconst apostrophe = \"\\u2019\";\nconst bullet1 =\n  \"Les d\u00e9fis rencontr\u00e9s c\" + apostrophe + \"est qu\" + apostrophe + \"au d\u00e9but je n\" + apostrophe +\n  \"arrivais pas \u00e0 supprim\u00e9 des fichiers donc j\" + apostrophe +\n  \"ai chercher une commande pour comment forcer la suppression.\";\nconst bullet2 =\n  \"J\" + apostrophe + \"ai un peu eu des difficult\u00e9s sur le code css mais je l\" + apostrophe +\n  \"ai r\u00e9gl\u00e9 tr\u00e8s rapidement\";\nconst bullet3 =\n  \"Au d\u00e9but je faisais git add * au lieu de git add . et je me disais pourquoi sa ne marchais pas et donc j\" +\n  apostrophe + \"ai su que c\" + apostrophe + \"\u00e9tait git add .\";\n\nconst lastPara = body.paragraphs.getLast();\nlastPara.load(\"text,style\");\nawait context.sync();\n\n// The three bullets are typed as one block of text, separated by paragraph marks\n// (\"\\r\"), so the existing empty ListParagraph/numId-6 paragraph receives the first\n// sentence and Word's own paragraph-split semantics mint the next two bullets with\n// the same style/numbering automatically (mirrors typing + Enter + Enter in the UI).\nlastPara.insertText(bullet1 + \"\\r\" + bullet2 + \"\\r\" + bullet3, Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n#\n# Reproduces the commit \"ajout des defis rencontrer sur le document word\":\n#   1. The empty paragraph right after the \"CREATION DE LA PAGE \u00ab SITE DORDI \u00bb \u00b8\"\n#      Heading 1 gets a single space \" \" typed into it.\n#   2. The (until now empty) last \"ListParagraph\" bullet (numId 6) at the end of\n#      the document is filled in with the first \"d\u00e9fis rencontr\u00e9s\" bullet, and two\n#      more bullets (same list/style) are appended after it.\n\n$d = $word.ActiveDocument\n\n# --- Locate paragraphs by scanning $d.Paragraphs (indices are reliable; the\n#     Paragraph.Next()/Previous() walk is not, around this document's TOC block) ---\n$count = $d.Paragraphs.Count\n\n# --- Change 1: the blank paragraph right after the \"CREATION DE LA PAGE...\" title ---\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$titleFound = $searchRange.Find.Execute(\"CREATION DE LA PAGE \u00ab SITE DORDI \u00bb \u00b8\")\n\n$titleIndex = -1\nif ($titleFound) {\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($searchRange.Start -ge $p.Range.Start -and $searchRange.Start -lt $p.Range.End) {\n            $titleIndex = $i\n            break\n        }\n    }\n}\n\nif ($titleIndex -gt 0 -and $titleIndex -lt $count) {\n    $blankPara = $d.Paragraphs.Item($titleIndex + 1)\n    $blankText = $blankPara.Range.Text -replace \"`r\", \"\"\n    if ($blankText -eq \"\") {\n        $blankPara.Range.InsertAfter(\" \")\n    }\n}\n\n# --- Change 2: fill in / append the \"d\u00e9fis rencontr\u00e9s\" bullet list ---\n# Find the LAST paragraph styled \"List Paragraph\" that is still empty - that is the\n# bullet placeholder waiting for the journal entry.\n$bulletIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Style.NameLocal\n    $txt = $p.Range.Text -replace \"`r\", \"\"\n    if ($styleName -eq \"List Paragraph\" -and $txt -eq \"\") {\n        $bulletIndex = $i\n    }\n}\n\nif ($bulletIndex -gt 0) {\n    $bulletPara = $d.Paragraphs.Item($bulletIndex)\n\n    $bullet1 = \"Les d\u00e9fis rencontr\u00e9s c\u2019est qu\u2019au d\u00e9but je n\u2019arrivais pas \u00e0 supprim\u00e9 des fichiers donc j\u2019ai chercher une commande pour comment forcer la suppression.\"\n    $bullet2 = \"J\u2019ai un peu eu des difficult\u00e9s sur le code css mais je l\u2019ai r\u00e9gl\u00e9 tr\u00e8s rapidement\"\n    $bullet3 = \"Au d\u00e9but je faisais git add * au lieu de git add . et je me disais pourquoi sa ne marchais pas et donc j\u2019ai su que c\u2019\u00e9tait git add .\"\n\n    # Typed as one block separated by paragraph marks (\"`r\"): the existing empty\n    # bullet receives the first sentence and Word mints the next two bullets with\n    # the same List Paragraph style / numId automatically (same as Enter + Enter).\n    $bulletPara.Range.InsertAfter($bullet1 + \"`r\" + $bullet2 + \"`r\" + $bullet3)\n}\n"}
